$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: same student id / name repeated with hyperlinks on column A
$email = "S530742@nwmissouri.edu"
$name = "Naveen "

for ($r = 2; $r -le 5; $r++) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 1), ("mailto:" + $email), "", "", $email)
    $ws.Cells.Item($r, 2).Value = $name
}

# Header row
$ws.Range("A1").Value = "SID"
$ws.Range("B1").Value = "Name"

$ws.Range("B2").Select()

$wb.Save()
